$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C12" = 0.158
    "C13" = -2.045
    "C39" = -0.53
    "C45" = -1.407
    "C57" = -1.84
    "C66" = -2.23
    "C76" = -2.807
    "C82" = -1.738
    "C90" = -0.797
    "C94" = 2.511
    "C97" = 1.12
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
